$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set Gulpease index value for AR column (G) on the "Periodi troppo lunghi o complessi" (row7)
# and "Parole non appropriate" (row8) rows to 0, and on "Indice gulpease" (row11) to 50.
$ws.Range("G7").Value = 0
$ws.Range("G8").Value = 0
$ws.Range("G11").Value = 50

# Update the active selection to G2 (matches the final cursor position in the diff)
$ws.Range("G2").Select()
